$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "VAV-L21-NW2"
$ws.Range("B6").Value = "L21-NW2"
$ws.Range("A7").Value = "VAV-L04-INT09"
$ws.Range("B7").Value = "L04_INT09"

$ws.Range("B8").Select()
